$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "de-poblacion-menor-de-25"
$ws.Range("E2").Value = "de-poblacion-menor-de-15"
$ws.Range("I2").Value = "de-poblacion-de-0-a-19-anos"
$ws.Range("K2").Value = "de-poblacion-de-65-y-mas-anos"
$ws.Range("Q2").Value = "de-poblacion-de-20-a-64-anos"
$ws.Range("R2").Value = "de-poblacion-menor-de-45"
$ws.Range("T2").Value = "de-poblacion-menor-de-35"

$ws.Range("D3").Value = "iaest-measure:de-poblacion-menor-de-25"
$ws.Range("E3").Value = "iaest-measure:de-poblacion-menor-de-15"
$ws.Range("I3").Value = "iaest-measure:de-poblacion-de-0-a-19-anos"
$ws.Range("K3").Value = "iaest-measure:de-poblacion-de-65-y-mas-anos"
$ws.Range("Q3").Value = "iaest-measure:de-poblacion-de-20-a-64-anos"
$ws.Range("R3").Value = "iaest-measure:de-poblacion-menor-de-45"
$ws.Range("T3").Value = "iaest-measure:de-poblacion-menor-de-35"
